$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.29699969291687
$ws.Range("B1").Value = 6.664901733398438
$ws.Range("C1").Value = 6.610248565673828
$ws.Range("D1").Value = 6.846599578857422
$ws.Range("E1").Value = 3.466582298278809
